$d = $word.ActiveDocument

# Locate the "שאלה 2" (Question 2) heading paragraph by its text, rather than
# assuming a fixed paragraph index.
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "שאלה 2") {
        $headingIndex = $i
        break
    }
}
if ($headingIndex -eq -1) {
    throw "Could not find the 'שאלה 2' heading paragraph"
}

# The paragraph right after the heading is the first existing "רצ""ב קוד"
# bullet (style "List Paragraph", numId 4, bidi). Insert a new paragraph
# immediately before it -- InsertParagraphBefore on that paragraph's range
# inherits its paragraph/run formatting -- then fill in the same text, so the
# new bullet matches its sibling list items under this question.
$firstCodePara = $d.Paragraphs.Item($headingIndex + 1)
$insertionPoint = $firstCodePara.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($headingIndex + 1)
$newRange = $newPara.Range.Duplicate
$newRange.MoveEnd(1, -1) | Out-Null
$newRange.Text = "רצ""ב קוד"
